# Update (Analyze PO & Forecast)
# - Shift all Week_Start_Date values in "Forecast Comparison" (B2:B17) back by 4 weeks (28 days)
# - Clear the is_holiday_week flags in column J (J2:J17) to blank
# - Shift the two derived week dates on the "Summary" sheet (B13, B15) back by 4 weeks (28 days) to match

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

function Set-WeekStartDate($ws, $row, $newValue) {
    $cell = $ws.Cells.Item($row, 2)
    # Force the value to be stored as literal text (not auto-converted to a date serial number)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    # Drop back to the default "Normal" style so no stray formatting is left behind
    $cell.Style = "Normal"
}

# Forecast Comparison!B2:B17 -> each date moves back 28 days, and J2:J17 is cleared to blank
Set-WeekStartDate $wsForecast 2  "2024-12-29"
Set-WeekStartDate $wsForecast 3  "2025-01-05"
Set-WeekStartDate $wsForecast 4  "2025-01-12"
Set-WeekStartDate $wsForecast 5  "2025-01-19"
Set-WeekStartDate $wsForecast 6  "2025-01-26"
Set-WeekStartDate $wsForecast 7  "2025-02-02"
Set-WeekStartDate $wsForecast 8  "2025-02-09"
Set-WeekStartDate $wsForecast 9  "2025-02-16"
Set-WeekStartDate $wsForecast 10 "2025-02-23"
Set-WeekStartDate $wsForecast 11 "2025-03-02"
Set-WeekStartDate $wsForecast 12 "2025-03-09"
Set-WeekStartDate $wsForecast 13 "2025-03-16"
Set-WeekStartDate $wsForecast 14 "2025-03-23"
Set-WeekStartDate $wsForecast 15 "2025-03-30"
Set-WeekStartDate $wsForecast 16 "2025-04-06"
Set-WeekStartDate $wsForecast 17 "2025-04-13"

for ($row = 2; $row -le 17; $row++) {
    $wsForecast.Cells.Item($row, 10).ClearContents()
}

# Summary sheet: the "Max Forecast Week" (B13) and "Min Forecast Week" (B15) values
# are derived from the same shifted weeks, so they move back 28 days too.
$wsSummary.Range("B13").NumberFormat = "@"
$wsSummary.Range("B13").Value = "2025-01-05"
$wsSummary.Range("B13").Style = "Normal"

$wsSummary.Range("B15").NumberFormat = "@"
$wsSummary.Range("B15").Value = "2025-01-12"
$wsSummary.Range("B15").Style = "Normal"
